# Add two new rows to the coverity scan sheet and update the filename text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$filename = "/home/rdkv-core/cov/cov-analysis-linux64-2023.6.0/bin/device/entservices-softwareupdate/MaintenanceManager/MaintenanceManager.cpp "

# Update the filename cell (C2) to include the new cov-analysis path prefix.
# This shared string is also referenced by C3 and C4 below.
$ws.Range("C2").Value = $filename

# Row 3: maintenanceManagerOnBootup / COPY_INSTEAD_OF_MOVE / filename / 1515
$ws.Range("A3").Value = "maintenanceManagerOnBootup"
$ws.Range("B3").Value = "COPY_INSTEAD_OF_MOVE"
$ws.Range("C3").Value = $filename
$ws.Range("D3").Value = 1515

# Row 4: setMaintenanceMode / COPY_INSTEAD_OF_MOVE / filename / 2288
$ws.Range("A4").Value = "setMaintenanceMode"
$ws.Range("B4").Value = "COPY_INSTEAD_OF_MOVE"
$ws.Range("C4").Value = $filename
$ws.Range("D4").Value = 2288

# Update the selected cell as reflected in the saved file.
$ws.Range("B6").Select()
